# Generate Report for Archive
#
# 1) The localization status "Ready for handoff" becomes "In Translation"
#    everywhere it is used (Overview sheet's zh-cn/de-de status columns,
#    and the Status column on each per-locale detail sheet).
# 2) The (now shorter) status text no longer needs as much horizontal
#    room, so the columns that display it are narrowed accordingly.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# New width (in "characters", the unit COM ColumnWidth uses) that best
# reproduces the narrower columns from the original workbook.
$newStatusColWidth = 12.5

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $used.Replace($oldStatus, $newStatus) | Out-Null
}

# Overview sheet: status columns are E (zh-cn) and F (de-de).
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = $newStatusColWidth
$overview.Columns.Item(6).ColumnWidth = $newStatusColWidth

# Per-locale detail sheets: status is column C.
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = $newStatusColWidth

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = $newStatusColWidth
